# Add a new "matches" worksheet in front of the existing sheets and fill it
# with the employee/position match data (added new data set).

$wb = $excel.ActiveWorkbook

# New sheets are inserted before the active sheet, i.e. at the front of the
# workbook - exactly where "matches" needs to land.
$ws = $wb.Worksheets.Add()
$ws.Name = "matches"

# Header row
$ws.Range("A1").Value = "employee_id"
$ws.Range("B1").Value = "position_id"
$ws.Range("C1").Value = "points"

# Match rows: employee_id / position_id / points
$data = @(
    @(1067, 1016, 100),
    @(1134, 1016, 500),
    @(1201, 1016, 500),
    @(1268, 1016, 500),
    @(1067, 1037, 1000),
    @(1134, 1037, 100),
    @(1201, 1037, 250),
    @(1268, 1037, 750),
    @(1067, 1033, 1000),
    @(1134, 1033, 250),
    @(1201, 1033, 0),
    @(1268, 1033, 0),
    @(1067, 1013, 100),
    @(1134, 1013, 250),
    @(1201, 1013, 400),
    @(1268, 1013, 100)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

# Match the author's final cursor position / page orientation on the new sheet
$ws.PageSetup.Orientation = 1
$null = $ws.Range("C2").Select()
